$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" column (G) values (replacing legacy Strike# data)
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
